$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.207.96"
$ws.Range("E2").Value = "  -2.26%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.369.36"
$ws.Range("E3").Value = "  -3.21%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "502.11"
$ws.Range("E5").Value = "  -1.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.36"
$ws.Range("E6").Value = "  -2.69%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.545"
$ws.Range("E8").Value = "  -2.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.369.86"
$ws.Range("E9").Value = "  -3.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0979"
$ws.Range("E10").Value = "  -0.12%  "

$ws.Range("E11").Value = "  +0.18%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.75"
$ws.Range("E12").Value = "  +3.36%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.323"
$ws.Range("E13").Value = "  +0.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.790.41"
$ws.Range("E14").Value = "  -3.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.169.81"
$ws.Range("E15").Value = "  -2.34%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.44"
$ws.Range("E16").Value = "  -1.78%  "

$ws.Range("E17").Value = "  -1.06%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.373.19"
$ws.Range("E18").Value = "  -1.97%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.00"
$ws.Range("E19").Value = "  -2.81%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.03"
$ws.Range("E20").Value = "  -1.89%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "306.72"
$ws.Range("E21").Value = "  -2.31%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.25"
$ws.Range("E22").Value = "  -2.41%  "

$ws.Range("E23").Value = "  +0.10%  "

$ws.Range("E24").Value = "  +0.26%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  +0.55%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.370"
$ws.Range("E26").Value = "  -3.00%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.146"
$ws.Range("E27").Value = "  -5.14%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.22"
$ws.Range("E28").Value = "  -4.32%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.92"
$ws.Range("E29").Value = "  -1.68%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0712"
$ws.Range("E30").Value = "  -2.85%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.64"
$ws.Range("E31").Value = "  -2.90%  "

$ws.Range("E32").Value = "  +0.24%  "

$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.08"
$ws.Range("E33").Value = "  -4.45%  "

$ws.Range("E34").Value = "  -6.89%  "

$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.997"
$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.55"
$ws.Range("E36").Value = "  -2.41%  "

$ws.Range("E37").Value = "  -5.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.75"
$ws.Range("E38").Value = "  -2.21%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.03"
$ws.Range("E39").Value = "  -1.69%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.791"
$ws.Range("E40").Value = "  -2.88%  "

$ws.Range("E41").Value = "  -5.57%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "129.90"
$ws.Range("E42").Value = "  -3.78%  "

$ws.Range("E43").Value = "  -1.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.72"
$ws.Range("E44").Value = "  -4.29%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.561"
$ws.Range("E45").Value = "  -1.93%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0904"
$ws.Range("E46").Value = "  -1.47%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "240.04"
$ws.Range("E47").Value = "  -6.10%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0480"
$ws.Range("E48").Value = "  -2.46%  "

$ws.Range("E49").Value = "  -2.90%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.03"
$ws.Range("E50").Value = "  +0.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.950"
$ws.Range("E51").Value = "  -0.75%  "
